# NIT-9002760823.xlsx - Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" lookup list backing column E (rows 16-56) is re-sorted
# from descending to ascending order (1802 .. 2106). Column E keeps the same
# row positions, so every label is rewritten in-place to the value that now
# falls at that slot once the list is sorted ascending. The "Valor Mora"
# (F) figures travel with their period label, so F16/F56 swap accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old (descending) order of the 41 period codes that filled E16:E56.
$oldOrder = @(
    "2106","2105","2104","2103","2102","2101",
    "2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802"
)

# New (ascending) order after the database refresh/sort.
$newOrder = @($oldOrder[($oldOrder.Length - 1)..0])

$firstRow = 16
for ($i = 0; $i -lt $oldOrder.Length; $i++) {
    $row = $firstRow + $i
    $newLabel = $newOrder[$i]
    $ws.Range("E$row").Value = $newLabel
}

# "Valor Mora" for the period that is now 2106 / 1802 swaps rows (16 <-> 56)
# along with the label re-sort above. Use Value2 (numeric) so the swap
# round-trips as a number rather than text.
$f16 = $ws.Range("F16").Value2
$f56 = $ws.Range("F56").Value2
$ws.Range("F16").Value2 = $f56
$ws.Range("F56").Value2 = $f16
